$d = $word.ActiveDocument

# 1) "We are using methods from" -> "We want to construct a food web for a paleo
#    community using two differnet methods; the pfim"
$d.Content.Find.Execute(
    "We are using methods from",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We want to construct a food web for a paleo community using two differnet methods; the pfim",
    1
)

# 2) "[1] and [2]" -> "[1] and the niche model [2]"
#    (scoped to this exact phrase so we don't touch the other "and"s
#     elsewhere in the document, e.g. in the caption / bibliography)
$d.Content.Find.Execute(
    "[1] and [2]",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[1] and the niche model [2]",
    1
)

# 3) "[2] to construct a food web for a paleo community. We want to see if the
#    different models are constructing different food webs" ->
#    "[2]. We want to see if the different models are constructing different
#    food webs"
$d.Content.Find.Execute(
    "[2] to construct a food web for a paleo community. We want to see if the different models are constructing different food webs",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[2]. We want to see if the different models are constructing different food webs",
    1
)
